# KeywordDictionaryForWebTesting.xlsx - "Updating repo with latestCodes"
#
# Adds a new keyword-dictionary row documenting a "textBoxShouldHaveValue"
# step (a Text Field function) right above the existing
# "waitForPageToRender" row, and switches the workbook back to automatic
# calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 17; this pushes the current row 17
# ("waitForPageToRender", Sno 16) down to row 18.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new keyword entry.
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Text Field"
$ws.Range("C17").Value = "textBoxShouldHaveValue"
$ws.Range("D17").Value = "Accepts two parameters @locator and @testData. It gets the text from textBox and validates against the @testData provided. If the validation fails testing should still continue"

# Match the wrapped-text formatting used by the rest of the table.
$ws.Range("C17:D17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 29

# The row that got shifted down to 18 keeps its original data but its
# Sno (column A) needs to be renumbered from 16 to 17.
$ws.Range("A18").Value = 17

# Mirror the author's final selection on the new row.
$ws.Range("C17").Select()

# The workbook was left in manual calculation mode; restore automatic.
$excel.Calculation = -4105  # xlCalculationAutomatic
